$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1981707317073171
$ws.Range("C2").Value = 0.5823170731707317
$ws.Range("J2").Value = 0.009146341463414634
$ws.Range("P2").Value = 0.1371951219512195
$ws.Range("S2").Value = 0.07317073170731707
$ws.Range("C3").Value = 0.01036269430051814
$ws.Range("J3").Value = 0.04663212435233161
$ws.Range("P3").Value = 0.7564766839378239
$ws.Range("S3").Value = 0.1865284974093264
$ws.Range("P4").Value = 0.6829268292682927
$ws.Range("S4").Value = 0.3170731707317073
$ws.Range("B6").Value = 0.065
$ws.Range("D6").Value = 0.01
$ws.Range("F6").Value = 0.02
$ws.Range("J6").Value = 0.24
$ws.Range("O6").Value = 0.015
$ws.Range("Q6").Value = 0.16
$ws.Range("R6").Value = 0.08500000000000001
$ws.Range("S6").Value = 0.405
$ws.Range("B7").Value = 0.1980676328502415
$ws.Range("D7").Value = 0.03864734299516908
$ws.Range("F7").Value = 0.01932367149758454
$ws.Range("J7").Value = 0.1256038647342995
$ws.Range("O7").Value = 0.004830917874396135
$ws.Range("Q7").Value = 0.1690821256038647
$ws.Range("R7").Value = 0.07246376811594203
$ws.Range("S7").Value = 0.3719806763285024
$ws.Range("B8").Value = 0.1113989637305699
$ws.Range("D8").Value = 0.01295336787564767
$ws.Range("F8").Value = 0.08031088082901554
$ws.Range("J8").Value = 0.1191709844559585
$ws.Range("O8").Value = 0.01036269430051814
$ws.Range("Q8").Value = 0.1761658031088083
$ws.Range("R8").Value = 0.08808290155440414
$ws.Range("S8").Value = 0.4015544041450777
$ws.Range("B9").Value = 0.1279620853080569
$ws.Range("D9").Value = 0.004739336492890996
$ws.Range("F9").Value = 0.06161137440758294
$ws.Range("J9").Value = 0.1279620853080569
$ws.Range("O9").Value = 0.02843601895734597
$ws.Range("Q9").Value = 0.1516587677725119
$ws.Range("R9").Value = 0.0995260663507109
$ws.Range("S9").Value = 0.3981042654028436
$ws.Range("B10").Value = 0.1053846153846154
$ws.Range("D10").Value = 0.02076923076923077
$ws.Range("E10").Value = 0.0007692307692307692
$ws.Range("F10").Value = 0.06615384615384616
$ws.Range("J10").Value = 0.1276923076923077
$ws.Range("O10").Value = 0.01
$ws.Range("Q10").Value = 0.2146153846153846
$ws.Range("R10").Value = 0.08307692307692308
$ws.Range("S10").Value = 0.3715384615384615
$ws.Range("G11").Value = 0.1442307692307692
$ws.Range("J11").Value = 0.08012820512820513
$ws.Range("K11").Value = 0.1987179487179487
$ws.Range("L11").Value = 0.5576923076923077
$ws.Range("S11").Value = 0.01923076923076923
$ws.Range("G12").Value = 0.7653631284916201
$ws.Range("J12").Value = 0.1843575418994413
$ws.Range("K12").Value = 0.00558659217877095
$ws.Range("L12").Value = 0.0335195530726257
$ws.Range("S12").Value = 0.0111731843575419
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.2916666666666667
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.02116402116402116
$ws.Range("H15").Value = 0.07407407407407407
$ws.Range("I15").Value = 0.08994708994708994
$ws.Range("J15").Value = 0.3915343915343915
$ws.Range("K15").Value = 0.08994708994708994
$ws.Range("M15").Value = 0.03174603174603174
$ws.Range("O15").Value = 0.04232804232804233
$ws.Range("S15").Value = 0.2592592592592592
$ws.Range("F16").Value = 0.009708737864077669
$ws.Range("H16").Value = 0.1359223300970874
$ws.Range("I16").Value = 0.0970873786407767
$ws.Range("J16").Value = 0.3932038834951456
$ws.Range("K16").Value = 0.1262135922330097
$ws.Range("M16").Value = 0.01941747572815534
$ws.Range("O16").Value = 0.06310679611650485
$ws.Range("S16").Value = 0.1553398058252427
$ws.Range("F17").Value = 0.01569506726457399
$ws.Range("H17").Value = 0.1726457399103139
$ws.Range("I17").Value = 0.09417040358744394
$ws.Range("J17").Value = 0.3968609865470852
$ws.Range("K17").Value = 0.1098654708520179
$ws.Range("M17").Value = 0.01345291479820628
$ws.Range("N17").Value = 0.002242152466367713
$ws.Range("O17").Value = 0.06726457399103139
$ws.Range("S17").Value = 0.1278026905829596
$ws.Range("F18").Value = 0.01030927835051546
$ws.Range("H18").Value = 0.1494845360824742
$ws.Range("I18").Value = 0.1030927835051546
$ws.Range("J18").Value = 0.4432989690721649
$ws.Range("K18").Value = 0.1134020618556701
$ws.Range("M18").Value = 0.0154639175257732
$ws.Range("O18").Value = 0.04123711340206185
$ws.Range("S18").Value = 0.1237113402061856
$ws.Range("F19").Value = 0.01402961808261886
$ws.Range("H19").Value = 0.186282151208106
$ws.Range("I19").Value = 0.08885424785658613
$ws.Range("J19").Value = 0.3865939204988308
$ws.Range("K19").Value = 0.102883865939205
$ws.Range("M19").Value = 0.02338269680436477
$ws.Range("N19").Value = 0.002338269680436477
$ws.Range("O19").Value = 0.06001558846453624
$ws.Range("S19").Value = 0.1356196414653157
